$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D and E columns (and B/C for rows 45-46) to Text format first,
# so Excel preserves values like "1.00" / "0.999" / percentages as literal text
# instead of auto-converting them to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '63.608.54'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '3.407.84'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '568.28'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").Value = '157.27'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.410.29'
$ws.Range("E8").Value = '  -0.56%  '
$ws.Range("D9").Value = '0.571'
$ws.Range("E9").Value = '  -7.89%  '
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("D11").Value = '0.119'
$ws.Range("E11").Value = '  -3.57%  '
$ws.Range("D12").Value = '0.423'
$ws.Range("E12").Value = '  -4.41%  '
$ws.Range("D13").Value = '3.991.44'
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '26.92'
$ws.Range("E15").Value = '  -3.84%  '
$ws.Range("D16").Value = '0.0000171'
$ws.Range("E16").Value = '  -9.38%  '
$ws.Range("D17").Value = '63.680.62'
$ws.Range("E17").Value = '  -1.53%  '
$ws.Range("D18").Value = '3.402.62'
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").Value = '6.09'
$ws.Range("E19").Value = '  -4.67%  '
$ws.Range("D20").Value = '13.55'
$ws.Range("E20").Value = '  -3.30%  '
$ws.Range("D21").Value = '383.50'
$ws.Range("E21").Value = '  +1.63%  '
$ws.Range("D22").Value = '7.76'
$ws.Range("E22").Value = '  -3.70%  '
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").Value = '70.95'
$ws.Range("E24").Value = '  -2.09%  '
$ws.Range("D25").Value = '0.515'
$ws.Range("E25").Value = '  -6.87%  '
$ws.Range("D26").Value = '0.0000114'
$ws.Range("E26").Value = '  -4.91%  '
$ws.Range("D27").Value = '9.68'
$ws.Range("E27").Value = '  -5.59%  '
$ws.Range("D28").Value = '0.178'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").Value = '6.06'
$ws.Range("E30").Value = '  -2.58%  '
$ws.Range("D31").Value = '1.39'
$ws.Range("E31").Value = '  -7.27%  '
$ws.Range("D32").Value = '1.98'
$ws.Range("E32").Value = '  -2.46%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '22.84'
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("D35").Value = '6.93'
$ws.Range("E35").Value = '  -4.25%  '
$ws.Range("D36").Value = '1.50'
$ws.Range("E36").Value = '  -6.62%  '
$ws.Range("D37").Value = '160.42'
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").Value = '0.840'
$ws.Range("E38").Value = '  +8.91%  '
$ws.Range("D39").Value = '1.82'
$ws.Range("E39").Value = '  -4.51%  '
$ws.Range("D40").Value = '2.821.29'
$ws.Range("E40").Value = '  -2.08%  '
$ws.Range("D41").Value = '25.91'
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("D42").Value = '42.98'
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").Value = '0.0717'
$ws.Range("E43").Value = '  -5.87%  '
$ws.Range("D44").Value = '6.36'
$ws.Range("E44").Value = '  -9.24%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '25.60'
$ws.Range("E45").Value = '  -3.98%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").Value = '4.35'
$ws.Range("E46").Value = '  -5.88%  '
$ws.Range("D47").Value = '0.0303'
$ws.Range("E47").Value = '  -3.80%  '
$ws.Range("D48").Value = '328.41'
$ws.Range("E48").Value = '  +2.04%  '
$ws.Range("D49").Value = '2.33'
$ws.Range("E49").Value = '  +6.98%  '
$ws.Range("D50").Value = '1.03'
$ws.Range("E50").Value = '  -4.85%  '
$ws.Range("D51").Value = '0.103'
$ws.Range("E51").Value = '  -5.68%  '
